$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Update the summary figures at the top of the statement
# ------------------------------------------------------------------
# VALOR MORA (total overdue amount) grows because of the two new rows
$ws.Range("E11").Value = 358984
# Cant. Periodos (number of overdue periods) goes from 2 to 3 (2507,2508,2509)
$ws.Range("F13").Value = 3

# ------------------------------------------------------------------
# 2) Make room for two new detail rows right after the current last
#    data row (row 21) - this naturally pushes the signature/footer
#    block (rows 26-27) down to rows 28-29, exactly like Excel does
#    when you insert rows through the UI.
# ------------------------------------------------------------------

# Remember the special "closing" border formatting that currently
# lives on row 21 (the last row of the table, with the thicker/
# final bottom border) so we can move it onto the new last row later.
$ws.Range("B21:J21").Copy()

# Insert two blank rows below row 21 (rows 22 and 23)
$ws.Rows("22:23").Insert()

# The new row 23 becomes the new last row of the table, so it gets
# the "closing" formatting that used to belong to row 21.
$ws.Range("B23:J23").PasteSpecial(-4122)

# Row 21 is no longer the last row, so it (and the other brand new
# row, 22) should look like a normal interior row - copy that look
# from row 20.
$ws.Range("B20:J20").Copy()
$ws.Range("B21:J22").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 3) Fill in the new period (2509) rows for the first two workers
# ------------------------------------------------------------------
$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "1002303693"
$ws.Range("D22").Value = "CLEINER RAFAEL PADILLA PADILLA"
$ws.Range("E22").Value = "2509"
$ws.Range("F22").Value = 56940
$ws.Range("G22").Value = 1423500

$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "1048936629"
$ws.Range("D23").Value = "DILSON RAFAEL CARO VARGAS"
$ws.Range("E23").Value = "2509"
$ws.Range("F23").Value = 60000
$ws.Range("G23").Value = 1500000
